$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:H1, matching style of existing headers (A1:E1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy formatting (bold font, border, centered alignment) from A1 to F1:H1
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Fill boolean (TRUE/FALSE) values for rows 2-21 in columns F, G, H
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 6).Value = $false
    $ws.Cells.Item($r, 7).Value = $false
    $ws.Cells.Item($r, 8).Value = $false
}

# Row 10 column F is TRUE per diff
$ws.Cells.Item(10, 6).Value = $true
